$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells are treated as text so values like
# "239.70" or "33.00" keep their exact string formatting instead of
# being auto-converted to numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.620.08"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "2.470.66"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "317.44"
$ws.Range("E5").Value = "  +1.27%  "
$ws.Range("D6").Value = "92.35"
$ws.Range("E6").Value = "  +0.02%  "
$ws.Range("E7").Value = "  +0.48%  "
$ws.Range("E8").Value = "  +0.22%  "
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").Value = "33.00"
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("E11").Value = "  +7.63%  "
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "2.850.24"
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "6.90"
$ws.Range("E14").Value = "  -0.56%  "
$ws.Range("D15").Value = "15.61"
$ws.Range("E15").Value = "  -4.99%  "
$ws.Range("D16").Value = "2.463.79"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("E17").Value = "  +1.47%  "
$ws.Range("D18").Value = "41.581.88"
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("D19").Value = "6.46"
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("E20").Value = "  +0.31%  "
$ws.Range("D21").Value = "71.17"
$ws.Range("E21").Value = "  -1.87%  "
$ws.Range("D22").Value = "11.33"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("D23").Value = "239.70"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "24.65"
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("E28").Value = "  +2.40%  "
$ws.Range("E29").Value = "  +1.56%  "
$ws.Range("D30").Value = "36.08"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "161.16"
$ws.Range("E31").Value = "  +2.27%  "
$ws.Range("D32").Value = "5.52"
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("E33").Value = "  +0.16%  "
$ws.Range("E34").Value = "  +0.35%  "
$ws.Range("D35").Value = "0.0766"
$ws.Range("E35").Value = "  +0.92%  "
$ws.Range("D36").Value = "17.26"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "1.85"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").Value = "2.91"
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("E39").Value = "  +1.56%  "
$ws.Range("D40").Value = "0.104"
$ws.Range("E40").Value = "  -2.44%  "
$ws.Range("D41").Value = "3.99"
$ws.Range("E41").Value = "  -1.67%  "
$ws.Range("D42").Value = "2.46"
$ws.Range("E42").Value = "  +2.63%  "
$ws.Range("D43").Value = "1.986.65"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "19.06"
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "0.0285"
$ws.Range("E45").Value = "  -0.17%  "
$ws.Range("D46").Value = "2.99"
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").Value = "9.19"
$ws.Range("E47").Value = "  +2.49%  "
$ws.Range("D48").Value = "2.706.32"
$ws.Range("E48").Value = "  -0.57%  "
$ws.Range("D49").Value = "97.42"
$ws.Range("E49").Value = "  -0.82%  "
$ws.Range("D50").Value = "74.19"
$ws.Range("E50").Value = "  +2.47%  "
$ws.Range("D51").Value = "67.28"
$ws.Range("E51").Value = "  -1.57%  "
